$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (Wins, Losses, Ties) after the existing
# "Unnamed: 28" column, matching the style of the other header cells.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record (Wins/Losses/Ties) for every player row.
for ($row = 2; $row -le 48; $row++) {
    $ws.Cells.Item($row, 30).Value = 56
    $ws.Cells.Item($row, 31).Value = 106
    $ws.Cells.Item($row, 32).Value = 0
}
